{"js": "// Correct output for {{FACHSTELLEN_KANTONAL}}:\n//  - Prefix the municipality bullet with a literal \"Gemeinde \" run.\n//  - Turn the old \"{{FACHSTELLEN_KANTONAL_LIST | multiline}}\" bullet into a\n//    proper Jinja for-loop over FACHSTELLEN_KANTONAL, matching the\n//    surrounding bullet list's numbering/style.\n\nconst paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their distinctive placeholder text.\nlet municipalityIdx = -1;\nlet fachstellenIdx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (\n    municipalityIdx === -1 &&\n    t.indexOf(\"{{MUNICIPALITY}}\") !== -1 &&\n    t.indexOf(\"Daten mit Adresse abmischen\") !== -1\n  ) {\n    municipalityIdx = i;\n  }\n  if (fachstellenIdx === -1 && t.indexOf(\"FACHSTELLEN_KANTONAL_LIST\") !== -1) {\n    fachstellenIdx = i;\n  }\n}\nif (municipalityIdx === -1) {\n  throw new Error(\"Could not find the {{MUNICIPALITY}} bullet paragraph\");\n}\nif (fachstellenIdx === -1) {\n  throw new Error(\"Could not find the FACHSTELLEN_KANTONAL_LIST bullet paragraph\");\n}\n\nconst municipalityPara = paras.items[municipalityIdx];\nconst fachstellenPara = paras.items[fachstellenIdx];\n\n// 1) Add a literal \"Gemeinde \" run in front of \"{{MUNICIPALITY}} \".\nmunicipalityPara.insertText(\"Gemeinde \", Word.InsertLocation.start);\n\n// 2) Append the Jinja loop-open tag after \"Daten mit Adresse abmischen\".\nmunicipalityPara.insertText(\n  \"{% for fachstelle in FACHSTELLEN_KANTONAL %}\",\n  Word.InsertLocation.end\n);\n\n// 3) Replace the old tabbed/no-bullet paragraph with a proper bullet-list\n//    paragraph (same style/numbering as its siblings) holding the loop body\n//    and the closing tag.\nconst flatOpcNs = \"http://schemas.microsoft.com/office/2006/xmlPackage\";\nconst partXml =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p>\" +\n  '<w:pPr><w:pStyle w:val=\"Aufzhlung\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">{{ fachstelle.NAME }} </w:t></w:r>' +\n  '<w:r><w:rPr><w:highlight w:val=\"green\"/></w:rPr><w:t>Daten mit Adresse abmischen</w:t></w:r>' +\n  \"<w:r><w:t>{% endfor %}</w:t></w:r>\" +\n  \"</w:p></w:body></w:document>\";\nconst pkg =\n  '<pkg:package xmlns:pkg=\"' +\n  flatOpcNs +\n  '\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  partXml +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nfachstellenPara.insertOoxml(pkg, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Correct output for {{FACHSTELLEN_KANTONAL}}:\n#  - Prefix the municipality bullet with a literal \"Gemeinde \" run.\n#  - Turn the old \"{{FACHSTELLEN_KANTONAL_LIST | multiline}}\" bullet into a\n#    proper Jinja for-loop over FACHSTELLEN_KANTONAL, matching the\n#    surrounding bullet list's numbering/style.\n\n$d = $word.ActiveDocument\n\n$municipalityPara = $null\n$fachstellenPara = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($municipalityPara -eq $null -and $t -like \"*MUNICIPALITY*\" -and $t -like \"*Daten mit Adresse abmischen*\") {\n        $municipalityPara = $p\n    }\n    if ($fachstellenPara -eq $null -and $t -like \"*FACHSTELLEN_KANTONAL_LIST*\") {\n        $fachstellenPara = $p\n    }\n}\n\nif ($municipalityPara -eq $null) {\n    throw \"Could not find the {{MUNICIPALITY}} bullet paragraph\"\n}\nif ($fachstellenPara -eq $null) {\n    throw \"Could not find the FACHSTELLEN_KANTONAL_LIST bullet paragraph\"\n}\n\n# 1) Add a literal \"Gemeinde \" run in front of \"{{MUNICIPALITY}} \".\n$municipalityPara.Range.InsertBefore(\"Gemeinde \")\n\n# 2) Append the Jinja loop-open tag after \"Daten mit Adresse abmischen\".\n$municipalityPara.Range.InsertAfter(\"{% for fachstelle in FACHSTELLEN_KANTONAL %}\")\n\n# 3) Replace the old tabbed/no-bullet paragraph with a proper bullet-list\n#    paragraph (same style/numbering as its siblings) holding the loop body\n#    and the closing tag.\n$partXml = '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Aufzhlung\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">{{ fachstelle.NAME }} </w:t></w:r><w:r><w:rPr><w:highlight w:val=\"green\"/></w:rPr><w:t>Daten mit Adresse abmischen</w:t></w:r><w:r><w:t>{% endfor %}</w:t></w:r></w:p></w:body></w:document>'\n$pkg = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + $partXml + '</pkg:xmlData></pkg:part></pkg:package>'\n\n$fachstellenPara.Range.InsertXML($pkg)\n"}
